$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: Police vs Sunrise
$ws.Range("F3").Value = "Police"
$ws.Range("G3").Value = 2
$ws.Range("H3").Value = "Sunrise"
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 1.72
$ws.Range("K3").Value = "19/08/2023 03:14"
$ws.Range("L3").Value = 1.72
$ws.Range("M3").Value = "19/08/2023 03:14"
$ws.Range("N3").Value = 3.15
$ws.Range("O3").Value = "19/08/2023 03:14"
$ws.Range("P3").Value = 3.25
$ws.Range("Q3").Value = "20/08/2023 13:03"
$ws.Range("R3").Value = 3.98
$ws.Range("S3").Value = "19/08/2023 03:14"
$ws.Range("T3").Value = 3.98
$ws.Range("U3").Value = "19/08/2023 03:14"
$ws.Range("V3").Value = "https://www.betexplorer.com/football/rwanda/premier-league/police-sunrise/8j5oa9Ep/"

# Row 4: Etincelles vs Gorilla
$ws.Range("F4").Value = "Etincelles"
$ws.Range("G4").Value = 1
$ws.Range("H4").Value = "Gorilla"
$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 2.17
$ws.Range("K4").Value = "19/08/2023 03:14"
$ws.Range("L4").Value = 2.02
$ws.Range("M4").Value = "19/08/2023 12:43"
$ws.Range("N4").Value = 2.94
$ws.Range("O4").Value = "19/08/2023 03:14"
$ws.Range("P4").Value = 3.03
$ws.Range("Q4").Value = "20/08/2023 13:06"
$ws.Range("R4").Value = 2.89
$ws.Range("S4").Value = "19/08/2023 03:14"
$ws.Range("T4").Value = 3.14
$ws.Range("U4").Value = "19/08/2023 12:43"
$ws.Range("V4").Value = "https://www.betexplorer.com/football/rwanda/premier-league/etincelles-gorilla/GfsIGypI/"

# Row 16: Etoile de L'Est vs APR
$ws.Range("F16").Value = "Etoile de L'Est"
$ws.Range("G16").Value = 0
$ws.Range("H16").Value = "APR"
$ws.Range("I16").Value = 1
$ws.Range("J16").Value = 6.85
$ws.Range("K16").Value = "02/09/2023 14:13"
$ws.Range("L16").Value = 6.49
$ws.Range("M16").Value = "02/09/2023 14:59"
$ws.Range("N16").Value = 4.52
$ws.Range("O16").Value = "02/09/2023 14:13"
$ws.Range("P16").Value = 3.85
$ws.Range("Q16").Value = "02/09/2023 14:59"
$ws.Range("R16").Value = 1.34
$ws.Range("S16").Value = "02/09/2023 14:13"
$ws.Range("T16").Value = 1.45
$ws.Range("U16").Value = "02/09/2023 14:58"
$ws.Range("V16").Value = "https://www.betexplorer.com/football/rwanda/premier-league/etoile-de-l-est-apr/Es3pl4Nt/"

# Row 17: Marines vs Etincelles
$ws.Range("F17").Value = "Marines"
$ws.Range("G17").Value = 1
$ws.Range("H17").Value = "Etincelles"
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 1.94
$ws.Range("K17").Value = "01/09/2023 03:13"
$ws.Range("L17").Value = 2.18
$ws.Range("M17").Value = "01/09/2023 04:34"
$ws.Range("N17").Value = 3.02
$ws.Range("O17").Value = "01/09/2023 03:13"
$ws.Range("P17").Value = 3.24
$ws.Range("Q17").Value = "02/09/2023 13:05"
$ws.Range("R17").Value = 3.31
$ws.Range("S17").Value = "01/09/2023 03:13"
$ws.Range("T17").Value = 2.92
$ws.Range("U17").Value = "01/09/2023 04:34"
$ws.Range("V17").Value = "https://www.betexplorer.com/football/rwanda/premier-league/marines-etincelles/d6hkmOxm/"

# Row 34: Mukura Victory Sports vs Sunrise
$ws.Range("F34").Value = "Mukura Victory Sports"
$ws.Range("G34").Value = 1
$ws.Range("H34").Value = "Sunrise"
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 1.95
$ws.Range("K34").Value = "09/10/2023 02:12"
$ws.Range("L34").Value = 1.64
$ws.Range("M34").Value = "10/10/2023 14:58"
$ws.Range("N34").Value = 2.94
$ws.Range("O34").Value = "09/10/2023 02:12"
$ws.Range("P34").Value = 3.36
$ws.Range("Q34").Value = "10/10/2023 14:58"
$ws.Range("R34").Value = 3.39
$ws.Range("S34").Value = "09/10/2023 02:12"
$ws.Range("T34").Value = 5.03
$ws.Range("U34").Value = "10/10/2023 14:57"
$ws.Range("V34").Value = "https://www.betexplorer.com/football/rwanda/premier-league/mukura-victory-sports-sunrise/AVRRFPZ4/"

# Row 35: Gorilla vs Gasogi United
$ws.Range("F35").Value = "Gorilla"
$ws.Range("G35").Value = 2
$ws.Range("H35").Value = "Gasogi United"
$ws.Range("I35").Value = 2
$ws.Range("J35").Value = 2.57
$ws.Range("K35").Value = "09/10/2023 02:12"
$ws.Range("L35").Value = 2.79
$ws.Range("M35").Value = "10/10/2023 14:49"
$ws.Range("N35").Value = 2.72
$ws.Range("O35").Value = "09/10/2023 02:12"
$ws.Range("P35").Value = 2.56
$ws.Range("Q35").Value = "10/10/2023 14:49"
$ws.Range("R35").Value = 2.57
$ws.Range("S35").Value = "09/10/2023 02:12"
$ws.Range("T35").Value = 2.81
$ws.Range("U35").Value = "10/10/2023 14:49"
$ws.Range("V35").Value = "https://www.betexplorer.com/football/rwanda/premier-league/gorilla-gasogi-united/25OJHo5h/"

# Row 50: Mukura Victory Sports vs Bugesera
$ws.Range("F50").Value = "Mukura Victory Sports"
$ws.Range("G50").Value = 2
$ws.Range("H50").Value = "Bugesera"
$ws.Range("I50").Value = 1
$ws.Range("J50").Value = 2.3
$ws.Range("K50").Value = "20/10/2023 02:12"
$ws.Range("L50").Value = 2.26
$ws.Range("M50").Value = "21/10/2023 14:38"
$ws.Range("N50").Value = 2.76
$ws.Range("O50").Value = "20/10/2023 02:12"
$ws.Range("P50").Value = 2.79
$ws.Range("Q50").Value = "21/10/2023 14:38"
$ws.Range("R50").Value = 2.86
$ws.Range("S50").Value = "20/10/2023 02:12"
$ws.Range("T50").Value = 3.29
$ws.Range("U50").Value = "21/10/2023 14:38"
$ws.Range("V50").Value = "https://www.betexplorer.com/football/rwanda/premier-league/mukura-victory-sports-bugesera/4OTMztUK/"

# Row 51: Muhazi United vs Musanze
$ws.Range("F51").Value = "Muhazi United"
$ws.Range("G51").Value = 1
$ws.Range("H51").Value = "Musanze"
$ws.Range("I51").Value = 1
$ws.Range("J51").Value = 2.91
$ws.Range("K51").Value = "20/10/2023 02:12"
$ws.Range("L51").Value = 3.09
$ws.Range("M51").Value = "21/10/2023 14:41"
$ws.Range("N51").Value = 2.71
$ws.Range("O51").Value = "20/10/2023 02:12"
$ws.Range("P51").Value = 2.7
$ws.Range("Q51").Value = "21/10/2023 14:41"
$ws.Range("R51").Value = 2.3
$ws.Range("S51").Value = "20/10/2023 02:12"
$ws.Range("T51").Value = 2.44
$ws.Range("U51").Value = "21/10/2023 14:41"
$ws.Range("V51").Value = "https://www.betexplorer.com/football/rwanda/premier-league/muhazi-united-musanze/0fFhVJir/"

# Row 62: Rayon Sport vs Mukura Victory Sports
$ws.Range("F62").Value = "Rayon Sport"
$ws.Range("G62").Value = 4
$ws.Range("H62").Value = "Mukura Victory Sports"
$ws.Range("I62").Value = 1
$ws.Range("J62").Value = 1.71
$ws.Range("K62").Value = "03/11/2023 02:13"
$ws.Range("L62").Value = 1.72
$ws.Range("M62").Value = "04/11/2023 11:35"
$ws.Range("N62").Value = 3.03
$ws.Range("O62").Value = "03/11/2023 02:13"
$ws.Range("P62").Value = 3.15
$ws.Range("Q62").Value = "04/11/2023 12:02"
$ws.Range("R62").Value = 4.24
$ws.Range("S62").Value = "03/11/2023 02:13"
$ws.Range("T62").Value = 4.75
$ws.Range("U62").Value = "04/11/2023 11:35"
$ws.Range("V62").Value = "https://www.betexplorer.com/football/rwanda/premier-league/rayon-sport-mukura-victory-sports/4pN1h5eN/"

# Row 63: Marines vs Amagaju
$ws.Range("F63").Value = "Marines"
$ws.Range("G63").Value = 1
$ws.Range("H63").Value = "Amagaju"
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 2.16
$ws.Range("K63").Value = "03/11/2023 02:13"
$ws.Range("L63").Value = 2.02
$ws.Range("M63").Value = "04/11/2023 13:04"
$ws.Range("N63").Value = 2.86
$ws.Range("O63").Value = "03/11/2023 02:13"
$ws.Range("P63").Value = 2.91
$ws.Range("Q63").Value = "04/11/2023 13:04"
$ws.Range("R63").Value = 2.99
$ws.Range("S63").Value = "03/11/2023 02:13"
$ws.Range("T63").Value = 3.76
$ws.Range("U63").Value = "04/11/2023 13:04"
$ws.Range("V63").Value = "https://www.betexplorer.com/football/rwanda/premier-league/marines-amagaju/ANNcgotH/"

# Row 64: Muhazi United vs APR
$ws.Range("F64").Value = "Muhazi United"
$ws.Range("G64").Value = 1
$ws.Range("H64").Value = "APR"
$ws.Range("I64").Value = 2
$ws.Range("J64").Value = 5.61
$ws.Range("K64").Value = "03/11/2023 02:13"
$ws.Range("L64").Value = 4.23
$ws.Range("M64").Value = "04/11/2023 13:55"
$ws.Range("N64").Value = 3.67
$ws.Range("O64").Value = "03/11/2023 02:13"
$ws.Range("P64").Value = 3
$ws.Range("Q64").Value = "04/11/2023 13:55"
$ws.Range("R64").Value = 1.43
$ws.Range("S64").Value = "03/11/2023 02:13"
$ws.Range("T64").Value = 1.87
$ws.Range("U64").Value = "04/11/2023 13:55"
$ws.Range("V64").Value = "https://www.betexplorer.com/football/rwanda/premier-league/muhazi-united-apr/ryM5iPAT/"

# Row 65: Musanze vs Kiyovu
$ws.Range("F65").Value = "Musanze"
$ws.Range("G65").Value = 1
$ws.Range("H65").Value = "Kiyovu"
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 2.91
$ws.Range("K65").Value = "03/11/2023 02:13"
$ws.Range("L65").Value = 2.12
$ws.Range("M65").Value = "04/11/2023 13:13"
$ws.Range("N65").Value = 2.71
$ws.Range("O65").Value = "03/11/2023 02:13"
$ws.Range("P65").Value = 3.05
$ws.Range("Q65").Value = "04/11/2023 13:49"
$ws.Range("R65").Value = 2.3
$ws.Range("S65").Value = "03/11/2023 02:13"
$ws.Range("T65").Value = 3.09
$ws.Range("U65").Value = "04/11/2023 13:13"
$ws.Range("V65").Value = "https://www.betexplorer.com/football/rwanda/premier-league/musanze-kiyovu/OU3YnqJp/"

# Row 66: new match AS Kigali vs Sunrise
$ws.Range("A66").Value = 65
$ws.Range("B66").Value = "rwanda"
$ws.Range("C66").Value = "premier-league"
$ws.Range("D66").Value = "2023-2024"
$ws.Range("E66").Value = 45235.58333333334
$ws.Range("F66").Value = "AS Kigali"
$ws.Range("G66").Value = 0
$ws.Range("H66").Value = "Sunrise"
$ws.Range("I66").Value = 1
$ws.Range("J66").Value = 1.63
$ws.Range("K66").Value = "04/11/2023 02:13"
$ws.Range("L66").Value = 1.3
$ws.Range("M66").Value = "05/11/2023 13:50"
$ws.Range("N66").Value = 3.1
$ws.Range("O66").Value = "04/11/2023 02:13"
$ws.Range("P66").Value = 4.11
$ws.Range("Q66").Value = "05/11/2023 13:50"
$ws.Range("R66").Value = 4.69
$ws.Range("S66").Value = "04/11/2023 02:13"
$ws.Range("T66").Value = 11.51
$ws.Range("U66").Value = "05/11/2023 13:50"
$ws.Range("V66").Value = "https://www.betexplorer.com/football/rwanda/premier-league/as-kigali-sunrise/GOjKAnBi/"

# Carry over the row-65 number formatting (bold/border index column + date style) to row 66
$ws.Range("A65").Copy()
$ws.Range("A66").PasteSpecial(-4122)
$ws.Range("E65").Copy()
$ws.Range("E66").PasteSpecial(-4122)
$excel.CutCopyMode = 0
